$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.020.46'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.384.00'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.77'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.13'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.66'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.122'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.04%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.963.88'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.85'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.386.23'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000170'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.109.91'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.11'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.67'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.95'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '384.44'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '76.51'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.14%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.57%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.94%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.48%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.22'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.95'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.30'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.94'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '165.58'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.419.10'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.99'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.10%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0765'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.52'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.80%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.776'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.84%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.461.25'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.84'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0262'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +9.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.205'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.03%  '
